# Applies the metadata update described by the commit diff for
# StructureDefinition-medication-item-provenance.xlsx:
#   1. Bump the "Date" property value.
#   2. Insert a new "Jurisdiction" property row (empty value) right after
#      "Contact", pushing "Description" and everything below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- 1. Update the "Date" property value (row 8, column B) ---------------
$ws.Cells.Item(8, 2).Value = "2024-09-12T14:01:50+00:00"

# --- 2. Insert a new "Jurisdiction" property row right after "Contact" ---
# "Contact" lives on row 10, so the new row goes in at row 11, pushing
# "Description" and everything below it down by one.
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
# A single apostrophe forces an (empty) text entry rather than a truly
# blank cell, matching the source workbook's "empty but string-typed"
# value cells (e.g. the many placeholder cells on the Elements sheet).
$ws.Cells.Item(11, 2).Value = "'"

# Re-apply the formatting of the row directly below (now row 12, the old
# "Description" row) onto the freshly inserted row so it matches the
# sheet's existing look (border/fill/alignment) instead of picking up
# whatever default/ad-hoc styling Insert()/text-coercion produced.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
